$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Duna RATE/s row (row 6): rebalance from placeholder "0C/s / ?C/s" to real values ---
$ws.Range("B6").Value = "0℃/s / -0,043℃/s"
$ws.Range("C6").Value = "0℃/s / -0,02℃/s"
$ws.Range("D6").Value = "0℃/s / -0,043℃/s"
$ws.Range("E6").Value = "0℃/s / -0,0616℃/s"

# --- Duna RATE/m row (row 17, left table): rebalance from placeholder "0C/m / ?C/m" to real values ---
$ws.Range("B17").Value = "0℃/m / -2,6℃/m"
$ws.Range("C17").Value = "0℃/m / -1,2℃/m"
$ws.Range("D17").Value = "0℃/m / -2,6℃/m"
$ws.Range("E17").Value = "0℃/m / -3,7℃/m"

# --- Duna TEMP TERRA row (row 17, right table): fill in previously empty cells ---
$ws.Range("H17").Value = "20℃/-15℃"
$ws.Range("I17").Value = "27℃/-8℃"
$ws.Range("J17").Value = "20℃/-15℃"
$ws.Range("K17").Value = "4℃/-23℃"

# --- Column width adjustments to fit the new longer strings ---
$ws.Columns.Item(2).ColumnWidth = 15.1666666666667
$ws.Columns.Item(3).ColumnWidth = 15.3333333333333
$ws.Columns.Item(4).ColumnWidth = 16.6666666666667
$ws.Columns.Item(5).ColumnWidth = 16.5

# --- Row height for row 6 (now wraps to two lines) ---
$ws.Rows.Item(6).RowHeight = 30

# --- Update view: scroll down and move selection ---
$excel.ActiveWindow.ScrollRow = 4
$ws.Range("H14").Select()
